# ============================================================
# Add 2022-Q3 data
#  1. Insert a new worksheet "2022-Q3" before the existing "2022-Q2" sheet
#     and populate it with the quarterly fund holdings data.
#  2. Update the "总计" (summary) sheet: insert a new row for 2022-Q3 at
#     the top of the data (row 2), pushing the older quarters down by one
#     row.
# ============================================================

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Part 1: insert the new "2022-Q3" worksheet ahead of "2022-Q2"
# ------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q3Sheet = $wb.Worksheets.Add($q2Sheet)
$q3Sheet.Name = "2022-Q3"

# Pull header / column-A formatting from the existing "2022-Q2" sheet so the
# new sheet matches the look of its siblings (bold header row + bordered
# index column) instead of plain default formatting.
$q2Sheet.Range("B1:H1").Copy()
$q3Sheet.Range("B1:H1").PasteSpecial(-4122)
$q2Sheet.Range("A2").Copy()
$q3Sheet.Range("A2:A36").PasteSpecial(-4122)

$q3Headers = @(
  "基金代码",
  "基金名称",
  "基金规模",
  "股票总仓位",
  "仓位占比",
  "持有市值(亿元)",
  "仓位排名"
)

for ($col = 0; $col -lt $q3Headers.Length; $col++) {
  $q3Sheet.Cells.Item(1, $col + 2).Value = $q3Headers[$col]
}

$q3Data = @(
  @("011162", "博时港股通领先趋势混合A", "15.21", "90.30", "8.04", "1.2229", 2),
  @("100061", "富国中国中小盘混合（QDII）人民币", "35.11", "83.32", "3.46", "1.2148", 5),
  @("010591", "富国中国中小盘混合（QDII）美元", "35.11", "83.32", "3.46", "1.2148", 5),
  @("005847", "富国沪港深业绩驱动混合A", "32.84", "85.47", "3.10", "1.0180", 10),
  @("007139", "富国民裕进取沪港深成长精选混合A", "10.56", "88.09", "9.17", "0.9684", 3),
  @("010714", "东方红远见价值混合", "15.24", "94.15", "5.58", "0.8504", 3),
  @("012588", "南方港股通优势企业混合A", "25.02", "81.50", "2.83", "0.7081", 9),
  @("007455", "富国蓝筹精选股票（QDII）人民币", "13.62", "85.59", "3.56", "0.4849", 8),
  @("010583", "富国蓝筹精选股票（QDII）美元", "13.62", "85.59", "3.56", "0.4849", 8),
  @("012208", "华夏港股前沿经济混合（QDII）A", "9.67", "89.48", "4.75", "0.4593", 5),
  @("010671", "景顺长城大中华混合（QDII）美元A", "9.42", "70.56", "4.67", "0.4399", 5),
  @("262001", "景顺长城大中华混合（QDII）人民币A", "9.42", "70.56", "4.67", "0.4399", 5),
  @("007368", "浙商沪港深精选混合A", "6.59", "84.00", "6.48", "0.4270", 3),
  @("012227", "景顺长城港股通全球竞争力混合A", "7.81", "74.96", "4.48", "0.3499", 6),
  @("011163", "博时港股通领先趋势混合C", "3.82", "90.30", "8.04", "0.3071", 2),
  @("011556", "富国民裕进取沪港深成长精选混合C", "2.46", "88.09", "9.17", "0.2256", 3),
  @("011635", "富国港股通策略精选混合A", "6.21", "73.36", "3.59", "0.2229", 8),
  @("100055", "富国全球科技互联网股票（QDII）", "3.95", "86.97", "4.16", "0.1643", 7),
  @("005228", "汇添富港股通专注成长混合", "4.41", "74.46", "3.38", "0.1491", 10),
  @("011117", "富国沪港深业绩驱动混合C", "3.53", "85.47", "3.10", "0.1094", 10),
  @("007182", "万家沪港深蓝筹混合A", "2.95", "87.92", "3.60", "0.1062", 9),
  @("013009", "万家港股通精选混合A", "1.78", "87.54", "3.45", "0.0614", 10),
  @("001215", "博时沪港深优质企业混合A", "2.22", "91.96", "2.71", "0.0602", 10),
  @("012589", "南方港股通优势企业混合C", "1.75", "81.50", "2.83", "0.0495", 9),
  @("012228", "景顺长城港股通全球竞争力混合C", "0.99", "74.96", "4.48", "0.0444", 6),
  @("006537", "恒生前海港股通精选混合", "0.95", "90.50", "4.42", "0.0420", 6),
  @("006781", "汇丰晋信港股通精选股票", "0.61", "89.76", "4.93", "0.0301", 1),
  @("013010", "万家港股通精选混合C", "0.71", "87.54", "3.45", "0.0245", 10),
  @("007369", "浙商沪港深精选混合C", "0.32", "84.00", "6.48", "0.0207", 3),
  @("011636", "富国港股通策略精选混合C", "0.57", "73.36", "3.59", "0.0205", 8),
  @("012209", "华夏港股前沿经济混合（QDII）C", "0.35", "89.48", "4.75", "0.0166", 5),
  @("007183", "万家沪港深蓝筹混合C", "0.35", "87.92", "3.60", "0.0126", 9),
  @("001942", "前海开源沪港深汇鑫灵活配置混合A", "0.17", "87.24", "4.80", "0.0082", 4),
  @("001943", "前海开源沪港深汇鑫灵活配置混合C", "0.09", "87.24", "4.80", "0.0043", 4),
  @("002555", "博时沪港深优质企业混合C", "0.08", "91.96", "2.71", "0.0022", 10)
)

for ($i = 0; $i -lt $q3Data.Length; $i++) {
  $row = $q3Data[$i]
  $r = $i + 2
  $q3Sheet.Cells.Item($r, 1).Value = $i
  $q3Sheet.Cells.Item($r, 2).Value = "'" + $row[0]
  $q3Sheet.Cells.Item($r, 3).Value = $row[1]
  $q3Sheet.Cells.Item($r, 4).Value = "'" + $row[2]
  $q3Sheet.Cells.Item($r, 5).Value = "'" + $row[3]
  $q3Sheet.Cells.Item($r, 6).Value = "'" + $row[4]
  $q3Sheet.Cells.Item($r, 7).Value = "'" + $row[5]
  $q3Sheet.Cells.Item($r, 8).Value = $row[6]
}

# ------------------------------------------------------------------
# Part 2: update the "总计" summary sheet - insert the 2022-Q3 totals as
# the new row 2 and shift the previously-existing rows down by one.
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Copy column-A formatting down into the newly-used row 6 before shifting
# values into it.
$summary.Cells.Item(5, 1).Copy()
$summary.Cells.Item(6, 1).PasteSpecial(-4122)
$summary.Cells.Item(6, 1).Value = 4

# Shift existing rows 2-5 down to rows 3-6 (processed bottom-up to avoid
# clobbering data before it is read).
$summary.Cells.Item(6, 2).Value = $summary.Cells.Item(5, 2).Value()
$summary.Cells.Item(6, 3).Value = $summary.Cells.Item(5, 3).Value()
$summary.Cells.Item(6, 4).Value = $summary.Cells.Item(5, 4).Value()

$summary.Cells.Item(5, 2).Value = $summary.Cells.Item(4, 2).Value()
$summary.Cells.Item(5, 3).Value = $summary.Cells.Item(4, 3).Value()
$summary.Cells.Item(5, 4).Value = $summary.Cells.Item(4, 4).Value()

$summary.Cells.Item(4, 2).Value = $summary.Cells.Item(3, 2).Value()
$summary.Cells.Item(4, 3).Value = $summary.Cells.Item(3, 3).Value()
$summary.Cells.Item(4, 4).Value = $summary.Cells.Item(3, 4).Value()

$summary.Cells.Item(3, 2).Value = $summary.Cells.Item(2, 2).Value()
$summary.Cells.Item(3, 3).Value = $summary.Cells.Item(2, 3).Value()
$summary.Cells.Item(3, 4).Value = $summary.Cells.Item(2, 4).Value()

# New row 2: the 2022-Q3 summary totals.
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 35
$summary.Cells.Item(2, 4).Value = 11.97
